$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310
$ws.Range("I33").Value = 318.46155
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 318.46155
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = -89.46154999999999
$ws.Range("N33").Value = -658

$ws.Range("H51").Value = 7666.6665
$ws.Range("I51").Value = 7500
$ws.Range("J51").Value = 8000
$ws.Range("K51").Value = 7500
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = -7016
$ws.Range("N51").Value = -8968

$ws.Range("H64").Value = 3868.0908
$ws.Range("I64").Value = 3166.3333
$ws.Range("J64").Value = 4131.25
$ws.Range("K64").Value = 3166.3333
$ws.Range("L64").Value = 4131.25
$ws.Range("M64").Value = -2918.3333
$ws.Range("N64").Value = -4627.25

$ws.Range("H67").Value = 3868.0908
$ws.Range("I67").Value = 3166.3333
$ws.Range("J67").Value = 4131.25
$ws.Range("K67").Value = 3166.3333
$ws.Range("L67").Value = 4131.25
$ws.Range("M67").Value = -2308.3333
$ws.Range("N67").Value = -5847.25

$ws.Range("H69").Value = 1650
$ws.Range("J69").Value = 1590.909
$ws.Range("L69").Value = 4772.727000000001
$ws.Range("N69").Value = -6520.727000000001

$ws.Range("H72").Value = 1650
$ws.Range("J72").Value = 1590.909
$ws.Range("L72").Value = 14318.181
$ws.Range("N72").Value = -23054.181

$ws.Range("H76").Value = 3388.6667
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3388.6667
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H112").Value = 1016.8125
$ws.Range("J112").Value = 1077.3334
$ws.Range("L112").Value = 3232.0002
$ws.Range("N112").Value = -5448.0002

$ws.Range("H116").Value = 4744
$ws.Range("I116").Value = 3400
$ws.Range("J116").Value = 5128
$ws.Range("K116").Value = 3400
$ws.Range("L116").Value = 5128
$ws.Range("M116").Value = 42
$ws.Range("N116").Value = -12012

$ws.Range("H135").Value = 26319444
$ws.Range("I135").Value = 1243.2858
$ws.Range("J135").Value = 100010410
$ws.Range("K135").Value = 11189.5722
$ws.Range("L135").Value = 900093690
$ws.Range("M135").Value = -8654.572200000001
$ws.Range("N135").Value = -900098760

$ws.Range("H137").Value = 1251.5667
$ws.Range("I137").Value = 1239.4584
$ws.Range("K137").Value = 3718.3752
$ws.Range("M137").Value = -1168.3752

$ws.Range("H138").Value = 28573848
$ws.Range("I138").Value = 58824988
$ws.Range("J138").Value = 3327.2778
$ws.Range("K138").Value = 176474964
$ws.Range("L138").Value = 9981.8334
$ws.Range("M138").Value = -176469824
$ws.Range("N138").Value = -20261.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 31073.75
$ws.Range("J119").Value = 31073.75
$ws.Range("L119").Value = 31073.75
$ws.Range("N119").Value = -40749.75

$ws.Range("H124").Value = 7300
$ws.Range("J124").Value = 7300
$ws.Range("L124").Value = 7300
$ws.Range("N124").Value = -17120

$ws.Range("H125").Value = 24860
$ws.Range("J125").Value = 24860
$ws.Range("L125").Value = 24860
$ws.Range("N125").Value = -34700

$ws.Range("H132").Value = 22533.459
$ws.Range("I132").Value = 1787.8125
$ws.Range("K132").Value = 5363.4375
$ws.Range("M132").Value = -2833.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1846.9565
$ws.Range("I105").Value = 1669.091
$ws.Range("J105").Value = 1902.8572
$ws.Range("K105").Value = 1669.091
$ws.Range("L105").Value = 1902.8572
$ws.Range("M105").Value = 77.90900000000011
$ws.Range("N105").Value = -5396.8572

$ws.Range("H107").Value = 1085.6
$ws.Range("I107").Value = 1228.75
$ws.Range("J107").Value = 513
$ws.Range("K107").Value = 1228.75
$ws.Range("L107").Value = 513
$ws.Range("M107").Value = 691.25
$ws.Range("N107").Value = -4353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 21715.92
$ws.Range("I58").Value = 1842.3636
$ws.Range("J58").Value = 37330.855
$ws.Range("K58").Value = 1842.3636
$ws.Range("L58").Value = 37330.855
$ws.Range("M58").Value = -1639.3636
$ws.Range("N58").Value = -37736.855

$ws.Range("H62").Value = 50003384
$ws.Range("I62").Value = 125003930
$ws.Range("J62").Value = 3022
$ws.Range("K62").Value = 125003930
$ws.Range("L62").Value = 3022
$ws.Range("M62").Value = -125003306
$ws.Range("N62").Value = -4270

$ws.Range("H65").Value = 50003384
$ws.Range("I65").Value = 125003930
$ws.Range("J65").Value = 3022
$ws.Range("K65").Value = 625019650
$ws.Range("L65").Value = 15110
$ws.Range("M65").Value = -625016530
$ws.Range("N65").Value = -21350

$ws.Range("H107").Value = 1845.2307
$ws.Range("I107").Value = 798.6
$ws.Range("J107").Value = 2499.375
$ws.Range("K107").Value = 798.6
$ws.Range("L107").Value = 2499.375
$ws.Range("M107").Value = 1121.4
$ws.Range("N107").Value = -6339.375

$ws.Range("H132").Value = 2788.087
$ws.Range("I132").Value = 1976.2106
$ws.Range("J132").Value = 6644.5
$ws.Range("K132").Value = 5928.6318
$ws.Range("L132").Value = 19933.5
$ws.Range("M132").Value = -3398.6318
$ws.Range("N132").Value = -24993.5

$ws.Range("H136").Value = 21715.92
$ws.Range("I136").Value = 1842.3636
$ws.Range("J136").Value = 37330.855
$ws.Range("K136").Value = 5527.0908
$ws.Range("L136").Value = 111992.565
$ws.Range("M136").Value = -2977.0908
$ws.Range("N136").Value = -117092.565

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 741.6667
$ws.Range("J45").Value = 975
$ws.Range("L45").Value = 2925
$ws.Range("N45").Value = -3989

$ws.Range("H117").Value = 1541.125
$ws.Range("J117").Value = 1480
$ws.Range("L117").Value = 4440
$ws.Range("N117").Value = -11324

$ws.Range("H130").Value = 2499.5
$ws.Range("I130").Value = 2000
$ws.Range("J130").Value = 2999
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 8997
$ws.Range("M130").Value = -980
$ws.Range("N130").Value = -19037

$ws.Range("H131").Value = 718.73627
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 718.73627
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2156.20881
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12236.20881

$ws.Range("H140").Value = 1799.8846
$ws.Range("I140").Value = 1093.4706
$ws.Range("J140").Value = 3134.2222
$ws.Range("K140").Value = 3280.4118
$ws.Range("L140").Value = 9402.6666
$ws.Range("M140").Value = 1899.5882
$ws.Range("N140").Value = -19762.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9341.25
$ws.Range("I113").Value = 14106.571
$ws.Range("J113").Value = 2669.8
$ws.Range("K113").Value = 14106.571
$ws.Range("L113").Value = 2669.8
$ws.Range("M113").Value = -11936.571
$ws.Range("N113").Value = -7009.8

$ws.Range("H132").Value = 40346.152
$ws.Range("I132").Value = 1104
$ws.Range("J132").Value = 103133.6
$ws.Range("K132").Value = 3312
$ws.Range("L132").Value = 309400.8
$ws.Range("M132").Value = -782
$ws.Range("N132").Value = -314460.8

$ws.Range("H134").Value = 24829.666
$ws.Range("J134").Value = 24829.666
$ws.Range("L134").Value = 74488.99800000001
$ws.Range("N134").Value = -79558.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2687.125
$ws.Range("J68").Value = 2739.4
$ws.Range("L68").Value = 2739.4
$ws.Range("N68").Value = -4237.4

$ws.Range("H71").Value = 2687.125
$ws.Range("J71").Value = 2739.4
$ws.Range("L71").Value = 13697
$ws.Range("N71").Value = -21185

$ws.Range("H122").Value = 703013.4399999999
$ws.Range("I122").Value = 935923.4399999999
$ws.Range("K122").Value = 2807770.32
$ws.Range("M122").Value = -2805320.32

$ws.Range("H127").Value = 36138.7
$ws.Range("J127").Value = 36138.7
$ws.Range("L127").Value = 36138.7
$ws.Range("N127").Value = -46058.7

$ws.Range("H132").Value = 1665.1428
$ws.Range("I132").Value = 1122.1538
$ws.Range("J132").Value = 2547.5
$ws.Range("K132").Value = 3366.4614
$ws.Range("L132").Value = 7642.5
$ws.Range("M132").Value = -836.4614000000001
$ws.Range("N132").Value = -12702.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 5300
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 5300
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -6548

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 5300
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 26500
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -32740

$ws.Range("H119").Value = 34999.5
$ws.Range("J119").Value = 34999.5
$ws.Range("L119").Value = 34999.5
$ws.Range("N119").Value = -44675.5

$ws.Range("H126").Value = 1396.9231
$ws.Range("J126").Value = 1747.5714
$ws.Range("L126").Value = 5242.7142
$ws.Range("N126").Value = -10182.7142

$ws.Range("H132").Value = 1107.5238
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060

$ws.Range("H136").Value = 25001488
$ws.Range("I136").Value = 32259206
$ws.Range("J136").Value = 2677.7778
$ws.Range("K136").Value = 96777618
$ws.Range("L136").Value = 8033.3334
$ws.Range("M136").Value = -96775068
$ws.Range("N136").Value = -13133.3334
